$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# --- ALC ---
$ws_ALC.Range("H6").Value = 100
$ws_ALC.Range("I6").Value = 100
$ws_ALC.Range("J6").Value = 100
$ws_ALC.Range("K6").Value = 300
$ws_ALC.Range("L6").Value = 300
$ws_ALC.Range("M6").Value = -188
$ws_ALC.Range("N6").Value = -524
$ws_ALC.Range("H100").Value = 3500
$ws_ALC.Range("I100").Value = 5000
$ws_ALC.Range("J100").Value = 2000
$ws_ALC.Range("K100").Value = 5000
$ws_ALC.Range("L100").Value = 2000
$ws_ALC.Range("M100").Value = -4459
$ws_ALC.Range("N100").Value = -3082
$ws_ALC.Range("H115").Value = 5666
$ws_ALC.Range("I115").Value = 5666
$ws_ALC.Range("K115").Value = 16998
$ws_ALC.Range("M115").Value = -15431
$ws_ALC.Range("H131").Value = 0
$ws_ALC.Range("I131").Value = 0
$ws_ALC.Range("K131").Value = 0
$ws_ALC.Range("M131").Value = $null
$ws_ALC.Range("H137").Value = 3691.3076
$ws_ALC.Range("I137").Value = 3365.2222
$ws_ALC.Range("J137").Value = 4425
$ws_ALC.Range("K137").Value = 10095.6666
$ws_ALC.Range("L137").Value = 13275
$ws_ALC.Range("M137").Value = -7545.6666
$ws_ALC.Range("N137").Value = -18375

# --- ARM ---
$ws_ARM.Range("H32").Value = 4111.933
$ws_ARM.Range("I32").Value = 4319.9287
$ws_ARM.Range("K32").Value = 4319.9287
$ws_ARM.Range("M32").Value = -4032.9287
$ws_ARM.Range("H110").Value = 1421.0625
$ws_ARM.Range("I110").Value = 1415.8
$ws_ARM.Range("K110").Value = 1415.8
$ws_ARM.Range("M110").Value = 629.2

# --- BSM ---
$ws_BSM.Range("H105").Value = 2751.1428
$ws_BSM.Range("I105").Value = 2020.6666
$ws_BSM.Range("J105").Value = 3299
$ws_BSM.Range("K105").Value = 2020.6666
$ws_BSM.Range("L105").Value = 3299
$ws_BSM.Range("M105").Value = -273.6666
$ws_BSM.Range("N105").Value = -6793
$ws_BSM.Range("H107").Value = 1970.3334
$ws_BSM.Range("I107").Value = 1455.5
$ws_BSM.Range("J107").Value = 3000
$ws_BSM.Range("K107").Value = 1455.5
$ws_BSM.Range("L107").Value = 3000
$ws_BSM.Range("M107").Value = 464.5
$ws_BSM.Range("N107").Value = -6840
$ws_BSM.Range("H139").Value = 40780
$ws_BSM.Range("J139").Value = 40780
$ws_BSM.Range("L139").Value = 40780
$ws_BSM.Range("N139").Value = -51060

# --- CRP ---
$ws_CRP.Range("H16").Value = 169452.17
$ws_CRP.Range("I16").Value = 252550
$ws_CRP.Range("J16").Value = 3256.5
$ws_CRP.Range("K16").Value = 252550
$ws_CRP.Range("L16").Value = 3256.5
$ws_CRP.Range("M16").Value = -252263
$ws_CRP.Range("N16").Value = -3830.5
$ws_CRP.Range("H31").Value = 2579.4
$ws_CRP.Range("I31").Value = 1489.2
$ws_CRP.Range("K31").Value = 1489.2
$ws_CRP.Range("M31").Value = -1194.2
$ws_CRP.Range("H34").Value = 2579.4
$ws_CRP.Range("I34").Value = 1489.2
$ws_CRP.Range("K34").Value = 1489.2
$ws_CRP.Range("M34").Value = -1287.2
$ws_CRP.Range("H105").Value = 3000
$ws_CRP.Range("I105").Value = 3000
$ws_CRP.Range("K105").Value = 3000
$ws_CRP.Range("M105").Value = -1253
$ws_CRP.Range("H107").Value = 112055.555
$ws_CRP.Range("I107").Value = 112055.555
$ws_CRP.Range("K107").Value = 112055.555
$ws_CRP.Range("M107").Value = -110135.555
$ws_CRP.Range("H113").Value = 169452.17
$ws_CRP.Range("I113").Value = 252550
$ws_CRP.Range("J113").Value = 3256.5
$ws_CRP.Range("K113").Value = 252550
$ws_CRP.Range("L113").Value = 3256.5
$ws_CRP.Range("M113").Value = -250380
$ws_CRP.Range("N113").Value = -7596.5

# --- CUL ---
$ws_CUL.Range("H56").Value = 15000
$ws_CUL.Range("I56").Value = 15000
$ws_CUL.Range("K56").Value = 15000
$ws_CUL.Range("M56").Value = -14470
$ws_CUL.Range("H59").Value = 1000
$ws_CUL.Range("I59").Value = 1000
$ws_CUL.Range("K59").Value = 3000
$ws_CUL.Range("M59").Value = -2460
$ws_CUL.Range("H68").Value = 998
$ws_CUL.Range("I68").Value = 0
$ws_CUL.Range("J68").Value = 998
$ws_CUL.Range("K68").Value = 0
$ws_CUL.Range("L68").Value = 2994
$ws_CUL.Range("M68").Value = $null
$ws_CUL.Range("N68").Value = -4616
$ws_CUL.Range("H71").Value = 998
$ws_CUL.Range("I71").Value = 0
$ws_CUL.Range("J71").Value = 998
$ws_CUL.Range("K71").Value = 0
$ws_CUL.Range("L71").Value = 8982
$ws_CUL.Range("M71").Value = $null
$ws_CUL.Range("N71").Value = -17094
$ws_CUL.Range("H80").Value = 7500
$ws_CUL.Range("I80").Value = 5250
$ws_CUL.Range("J80").Value = 12000
$ws_CUL.Range("K80").Value = 15750
$ws_CUL.Range("L80").Value = 36000
$ws_CUL.Range("M80").Value = -14814
$ws_CUL.Range("N80").Value = -37872
$ws_CUL.Range("H83").Value = 7500
$ws_CUL.Range("I83").Value = 5250
$ws_CUL.Range("J83").Value = 12000
$ws_CUL.Range("K83").Value = 47250
$ws_CUL.Range("L83").Value = 108000
$ws_CUL.Range("M83").Value = -42570
$ws_CUL.Range("N83").Value = -117360
$ws_CUL.Range("H97").Value = 704.5
$ws_CUL.Range("I97").Value = 606
$ws_CUL.Range("J97").Value = 1000
$ws_CUL.Range("K97").Value = 1818
$ws_CUL.Range("L97").Value = 3000
$ws_CUL.Range("M97").Value = -1322
$ws_CUL.Range("N97").Value = -3992
$ws_CUL.Range("H131").Value = 919.8333
$ws_CUL.Range("I131").Value = 803.8
$ws_CUL.Range("K131").Value = 2411.4
$ws_CUL.Range("M131").Value = 2628.6
$ws_CUL.Range("H137").Value = 6300
$ws_CUL.Range("J137").Value = 0
$ws_CUL.Range("L137").Value = 0
$ws_CUL.Range("N137").Value = $null

# --- LTW ---
$ws_LTW.Range("H61").Value = 2957.5715
$ws_LTW.Range("I61").Value = 2740.8
$ws_LTW.Range("J61").Value = 3499.5
$ws_LTW.Range("K61").Value = 2740.8
$ws_LTW.Range("L61").Value = 3499.5
$ws_LTW.Range("M61").Value = -2538.8
$ws_LTW.Range("N61").Value = -3903.5
$ws_LTW.Range("H93").Value = 1700
$ws_LTW.Range("I93").Value = 1700
$ws_LTW.Range("K93").Value = 1700
$ws_LTW.Range("M93").Value = -452
$ws_LTW.Range("H111").Value = 99387
$ws_LTW.Range("J111").Value = 99387
$ws_LTW.Range("L111").Value = 99387
$ws_LTW.Range("N111").Value = -107567
$ws_LTW.Range("H113").Value = 2957.5715
$ws_LTW.Range("I113").Value = 2740.8
$ws_LTW.Range("J113").Value = 3499.5
$ws_LTW.Range("K113").Value = 2740.8
$ws_LTW.Range("L113").Value = 3499.5
$ws_LTW.Range("M113").Value = -570.8000000000002
$ws_LTW.Range("N113").Value = -7839.5
$ws_LTW.Range("H122").Value = 0
$ws_LTW.Range("I122").Value = 0
$ws_LTW.Range("K122").Value = 0
$ws_LTW.Range("M122").Value = $null

# --- WVR ---
$ws_WVR.Range("H93").Value = 0
$ws_WVR.Range("I93").Value = 0
$ws_WVR.Range("K93").Value = 0
$ws_WVR.Range("M93").Value = $null
$ws_WVR.Range("H132").Value = 4173
$ws_WVR.Range("I132").Value = 3459.6667
$ws_WVR.Range("K132").Value = 10379.0001
$ws_WVR.Range("M132").Value = -7849.000100000001

